$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLeft alignment constant
$xlLeft = -4131

# New device rows (157-161), following the exact pattern of the preceding
# rows: id, name, mac_address, serial_num, dspec_id, lang_code(eng),
# is_active(TRUE, left-aligned), cr_by(superadmin), cr_dtimes(now()),
# eff_dtimes(now()).
$newRows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; Dspec = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";         Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; Dspec = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";           Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; Dspec = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";     Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; Dspec = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";               Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; Dspec = 920 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Id
    $ws.Cells.Item($rowNum, 2).Value = $r.Name
    $ws.Cells.Item($rowNum, 3).Value = $r.Mac
    $ws.Cells.Item($rowNum, 4).Value = $r.Serial
    $ws.Cells.Item($rowNum, 6).Value = $r.Dspec
    $ws.Cells.Item($rowNum, 7).Value = "eng"
    $ws.Cells.Item($rowNum, 8).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($rowNum, 8).Value = $true
    $ws.Cells.Item($rowNum, 9).Value = "superadmin"
    $ws.Cells.Item($rowNum, 10).Value = "now()"
    $ws.Cells.Item($rowNum, 11).Value = "now()"
}

# Five trailing blank rows (162-166) that only carry the left-aligned
# style previously used for column H (is_active), same as the sheet's
# original ragged tail formatting.
for ($rowNum = 162; $rowNum -le 166; $rowNum++) {
    $ws.Cells.Item($rowNum, 8).HorizontalAlignment = $xlLeft
}

# Move the selection/view down to the newly added data, mirroring the
# author's on-screen position after typing the new rows.
[void]$ws.Range("E159").Select()
